$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook

# Sheet 1: "Сравнительные характеристики" - add two new model rows (7 and 8)
$ws1 = $wb.Worksheets.Item(1)

# Row 7: FFF model - dominated/equal case, D column repeats the "23.001" text value
$ws1.Range("A7").Value = "FFF"
$ws1.Range("B7").Value = 23
$ws1.Range("C7").Value = 24
# Copy D6 (already stored as text "23.001") into D7 so it keeps the same text type
$ws1.Range("D6").Copy()
$ws1.Range("D7").PasteSpecial($xlPasteValues)

# Row 8: GGG model - plain numeric row
$ws1.Range("A8").Value = "GGG"
$ws1.Range("B8").Value = 4
$ws1.Range("C8").Value = 55
$ws1.Range("D8").Value = 12

# Move the active selection on sheet 1 to F2
$ws1.Select()
$ws1.Range("F2").Select()

# Sheet 2: "System" - change criteria count (B1) from 5 to 7
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = 7
